$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new visitor row (row 4) mirroring the existing row 3 entry,
# with a fresh timestamp and userAgent (new VisitTracker 'viewedAtLocal'
# style entry appended by the API), reusing ip / pathname / referrer / language.
$ws.Range("A4").Value = "2026-01-21T07:34:37.311Z"
$ws.Range("B4").Value = "::1"
$ws.Range("C4").Value = "Mozilla/5.0 (Windows NT 10.0; Win64; x64) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/144.0.0.0 Safari/537.36"
$ws.Range("D4").Value = "/"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "en-US"

# Restore page-setup print properties that accompany the refreshed export.
$ws.PageSetup.FirstPageNumber = 1
$ws.PageSetup.UseFirstPageNumber = $true
$ws.PageSetup.Copies = 1
